{"js": "// Update the date heading and the 20 division problems in the practice\n// table. Cells are addressed by their (row, column) position in the\n// table so the edit is unambiguous even though several of the new\n// values duplicate other old values elsewhere in the document.\n\nconst body = context.document.body;\n\n// --- Title paragraph: date update -----------------------------------\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\nparagraphs.items[0].insertText(\"2025-09-14 Sunday\", Word.InsertLocation.replace);\n\n// --- Table of division problems --------------------------------------\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Rows that actually contain problems (0-based): rows 0, 4, 8, 12, 16 each\n// hold 5 filled cells; the rows in between are blank \"answer\" rows.\nconst dataRows = [0, 4, 8, 12, 16];\nconst newValues = [\n  [\"51\u00f76=\", \"40\u00f76=\", \"10\u00f74=\", \"47\u00f77=\", \"61\u00f76=\"],\n  [\"73\u00f76=\", \"30\u00f74=\", \"57\u00f79=\", \"29\u00f76=\", \"63\u00f78=\"],\n  [\"88\u00f76=\", \"53\u00f74=\", \"37\u00f78=\", \"43\u00f74=\", \"25\u00f74=\"],\n  [\"37\u00f76=\", \"72\u00f72=\", \"76\u00f79=\", \"22\u00f73=\", \"17\u00f76=\"],\n  [\"11\u00f72=\", \"96\u00f72=\", \"85\u00f76=\", \"75\u00f79=\", \"29\u00f78=\"],\n];\n\nfor (let i = 0; i < dataRows.length; i++) {\n  const rowIndex = dataRows[i];\n  const rowValues = newValues[i];\n  for (let col = 0; col < rowValues.length; col++) {\n    const cell = table.getCell(rowIndex, col);\n    cell.value = rowValues[col];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date heading and the 20 division problems in the practice\n# table. Cells are addressed by their (row, column) position in the table\n# so that the edit is unambiguous even though several of the new values\n# duplicate other old values elsewhere in the document.\n\n$d = $word.ActiveDocument\n\n# --- Title paragraph: date update -----------------------------------\n$d.Paragraphs.Item(1).Range.Text = \"2025-09-14 Sunday\"\n\n# --- Table of division problems --------------------------------------\n$t = $d.Tables.Item(1)\n\n# Rows that actually contain problems (1-based, matching the COM model):\n# rows 1, 5, 9, 13, 17 each hold 5 filled cells; the rows in between are\n# blank \"answer\" rows.\n$newValues = @(\n    @(\"51\u00f76=\", \"40\u00f76=\", \"10\u00f74=\", \"47\u00f77=\", \"61\u00f76=\"),\n    @(\"73\u00f76=\", \"30\u00f74=\", \"57\u00f79=\", \"29\u00f76=\", \"63\u00f78=\"),\n    @(\"88\u00f76=\", \"53\u00f74=\", \"37\u00f78=\", \"43\u00f74=\", \"25\u00f74=\"),\n    @(\"37\u00f76=\", \"72\u00f72=\", \"76\u00f79=\", \"22\u00f73=\", \"17\u00f76=\"),\n    @(\"11\u00f72=\", \"96\u00f72=\", \"85\u00f76=\", \"75\u00f79=\", \"29\u00f78=\")\n)\n\n$dataRows = @(1, 5, 9, 13, 17)\n\nfor ($i = 0; $i -lt $dataRows.Count; $i++) {\n    $rowIndex = $dataRows[$i]\n    $rowValues = $newValues[$i]\n    for ($col = 1; $col -le 5; $col++) {\n        $cell = $t.Cell($rowIndex, $col)\n        $cell.Range.Text = $rowValues[$col - 1]\n    }\n}\n"}
